$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '327.34'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '0.96%'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '44.00'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '-1.22%'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.495'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '-0.10%'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.08012'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '-0.26%'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '2.015'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '7.23%'
$ws.Range('B7').Value = 'GateToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '4.317'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '0.02%'
$ws.Range('B8').Value = 'BTSEToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '2.569'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '-3.16%'
$ws.Range('B9').Value = 'MXToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.9485'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '0.95%'
$ws.Range('B10').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C10').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.1122'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '-4.23%'
$ws.Range('B11').Value = 'WazirX'
$ws.Range('C11').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.1858'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '-0.97%'
$ws.Range('B12').Value = 'MCDex'
$ws.Range('C12').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '10.62'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '25.30%'
$ws.Range('B13').Value = 'MandalaExchangeToken'
$ws.Range('C13').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.09915'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '0.13%'
$ws.Range('B14').Value = 'BitrueCoin'
$ws.Range('C14').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.04592'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '9.39%'
$ws.Range('B15').Value = 'BitMartToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.1066'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '0.12%'
$ws.Range('B16').Value = 'BitForexToken'
$ws.Range('C16').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.001277'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '-0.08%'
$ws.Range('B17').Value = 'CoinExToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.04073'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '-4.28%'
$ws.Range('B18').Value = 'TigerCash'
$ws.Range('C18').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.005924'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '-1.43%'
$ws.Range('B19').Value = 'LEO'
$ws.Range('C19').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.354'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '-6.62%'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.3475'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '-0.31%'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.1409'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '2.49%'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '-3.84%'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '1.38%'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.004322'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '-3.27%'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '-6.17%'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0003739'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '-6.66%'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02572'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '-2.49%'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05680'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '3.73%'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.007542'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '-1.47%'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.007590'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '9.09%'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '-2.04%'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.008378'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '-8.82%'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00007093'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '-0.53%'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '-0.56%'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '54.81%'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.003546'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '3.14%'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.00002097'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '-0.56%'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0001997'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '-0.56%'
